$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Coin/Link/Price/Volume columns to text format before writing,
# so that numeric-looking values (e.g. "1.00", "8.30") are stored verbatim
# as strings instead of being coerced into numbers by Excel.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '68.521.57'
$ws.Range("E2").Value = '  -0.80%  '

$ws.Range("D3").Value = '3.473.11'
$ws.Range("E3").Value = '  -1.35%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '590.64'
$ws.Range("E5").Value = '  +2.11%  '

$ws.Range("D6").Value = '167.56'
$ws.Range("E6").Value = '  -2.16%  '

$ws.Range("D7").Value = '0.606'
$ws.Range("E7").Value = '  -2.33%  '

$ws.Range("D8").Value = '3.466.32'
$ws.Range("E8").Value = '  -1.18%  '

$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("E10").Value = '  +0.78%  '

$ws.Range("D11").Value = '6.79'
$ws.Range("E11").Value = '  +1.21%  '

$ws.Range("D12").Value = '0.571'
$ws.Range("E12").Value = '  -4.84%  '

$ws.Range("D13").Value = '46.49'
$ws.Range("E13").Value = '  -1.68%  '

$ws.Range("E14").Value = '  +0.86%  '

$ws.Range("D15").Value = '4.038.13'
$ws.Range("E15").Value = '  -1.10%  '

$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Value = '8.30'
$ws.Range("E16").Value = '  -5.63%  '

$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").Value = '612.14'
$ws.Range("E17").Value = '  -10.67%  '

$ws.Range("D18").Value = '3.483.18'
$ws.Range("E18").Value = '  -0.86%  '

$ws.Range("D19").Value = '68.693.07'
$ws.Range("E19").Value = '  -0.67%  '

$ws.Range("E20").Value = '  -2.16%  '

$ws.Range("D21").Value = '17.16'
$ws.Range("E21").Value = '  -1.47%  '

$ws.Range("D22").Value = '11.08'
$ws.Range("E22").Value = '  -0.74%  '

$ws.Range("D23").Value = '0.869'
$ws.Range("E23").Value = '  -4.18%  '

$ws.Range("D24").Value = '15.74'
$ws.Range("E24").Value = '  -4.97%  '

$ws.Range("D25").Value = '95.58'
$ws.Range("E25").Value = '  -2.12%  '

$ws.Range("E26").Value = '  -1.39%  '

$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = '5.82'
$ws.Range("E27").Value = '  +1.45%  '

$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.02%  '

$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").Value = '2.60'
$ws.Range("E29").Value = '  -2.21%  '

$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '9.06'
$ws.Range("E30").Value = '  -3.72%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '32.77'
$ws.Range("E31").Value = '  -1.53%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '8.38'
$ws.Range("E32").Value = '  -5.14%  '

$ws.Range("B33").Value = 'Stacks'
$ws.Range("C33").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D33").Value = '3.07'
$ws.Range("E33").Value = '  -3.39%  '

$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D34").Value = '1.31'
$ws.Range("E34").Value = '  -3.20%  '

$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").Value = '6.77'
$ws.Range("E35").Value = '  -6.54%  '

$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").Value = '571.07'
$ws.Range("E36").Value = '  -0.15%  '

$ws.Range("B37").Value = 'Cosmos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D37").Value = '10.68'
$ws.Range("E37").Value = '  -1.52%  '

$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").Value = '3.48'
$ws.Range("E38").Value = '  -4.65%  '

$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = '56.90'
$ws.Range("E39").Value = '  -0.40%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '0.101'
$ws.Range("E40").Value = '  -4.25%  '

$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.06%  '

$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '0.137'
$ws.Range("E42").Value = '  -0.89%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '0.0437'
$ws.Range("E43").Value = '  -0.50%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '3.387.06'
$ws.Range("E44").Value = '  -1.60%  '

$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").Value = '0.322'
$ws.Range("E45").Value = '  -4.43%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '32.42'
$ws.Range("E46").Value = '  -2.46%  '

$ws.Range("B47").Value = 'PEPE'
$ws.Range("C47").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D47").Value = '0.0₃0691'
$ws.Range("E47").Value = '  -1.67%  '

$ws.Range("B48").Value = 'ThetaToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D48").Value = '2.83'
$ws.Range("E48").Value = '  -1.61%  '

$ws.Range("B49").Value = 'Fetch.AI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D49").Value = '2.54'
$ws.Range("E49").Value = '  -1.62%  '

$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = '0.128'
$ws.Range("E50").Value = '  -4.07%  '

$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = '132.40'
$ws.Range("E51").Value = '  -1.34%  '

